$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.082.15"
$ws.Range("E2").Value = "  +2.26%  "

$ws.Range("D3").Value = "2.349.80"
$ws.Range("E3").Value = "  +7.04%  "

$ws.Range("E4").Value = "  -0.90%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.71"
$ws.Range("E5").Value = "  +5.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.66"
$ws.Range("E6").Value = "  +2.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.638"
$ws.Range("E7").Value = "  +3.51%  "

$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("E9").Value = "  +6.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.25"
$ws.Range("E10").Value = "  -2.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("E11").Value = "  +4.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.85"
$ws.Range("E12").Value = "  +2.90%  "

$ws.Range("E13").Value = "  +14.70%  "

$ws.Range("E14").Value = "  +2.07%  "

$ws.Range("E15").Value = "  +11.45%  "

$ws.Range("D16").Value = "2.715.55"
$ws.Range("E16").Value = "  +7.17%  "

$ws.Range("D17").Value = "2.451.36"
$ws.Range("E17").Value = "  +9.18%  "

$ws.Range("D18").Value = "43.075.52"
$ws.Range("E18").Value = "  +2.47%  "

$ws.Range("E19").Value = "  +4.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  +2.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.23"
$ws.Range("E21").Value = "  +3.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.40"
$ws.Range("E22").Value = "  +1.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.51"
$ws.Range("E23").Value = "  +12.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "252.05"
$ws.Range("E24").Value = "  +12.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.92"
$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("E26").Value = "  +4.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.96"
$ws.Range("E28").Value = "  +4.39%  "

$ws.Range("E29").Value = "  +0.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.50"
$ws.Range("E30").Value = "  +9.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.06"
$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("E32").Value = "  -0.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0920"
$ws.Range("E33").Value = "  +7.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.94"
$ws.Range("E34").Value = "  +8.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.02"
$ws.Range("E35").Value = "  +6.09%  "

$ws.Range("E36").Value = "  +6.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0377"
$ws.Range("E37").Value = "  +6.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.07"
$ws.Range("E38").Value = "  -1.20%  "

$ws.Range("E39").Value = "  +1.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  +12.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.49"
$ws.Range("E41").Value = "  +17.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.86"
$ws.Range("E42").Value = "  +4.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.230"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.36"
$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.60"
$ws.Range("E46").Value = "  +4.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.33"
$ws.Range("E47").Value = "  +11.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.82"
$ws.Range("E48").Value = "  +8.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("E49").Value = "  +1.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0999"
$ws.Range("E50").Value = "  +2.29%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "69.88"
$ws.Range("E51").Value = "  +9.33%  "
